# ---------------------------------------------------------------------------
# [ADDITIONAL SCRAPING] added scraping code for extra bowling attributes and
# excel sheets.
#
# 1) "ODI Batting Extra" had a bunch of cells (columns B-E) that only ever
#    held an empty string -- strip those down to genuinely-blank cells.
# 2) Add a brand new "ODI Bowling Extra" sheet (after "ODI Batting Extra")
#    with MATCH_CODE / MAIDEN_OVERS / PERCENT_WICKETS_OF_ALL columns.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Part 1: drop the empty placeholder cells left over on "ODI Batting Extra"
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")

$emptyCells = @(
    "B3","C3","D3","E3",
    "E4",
    "E6",
    "C8","D8","E8",
    "B10","C10","D10","E10",
    "C12","D12","E12",
    "B13","C13","D13","E13",
    "B17","C17","D17","E17",
    "B18","C18","D18","E18",
    "B20","C20","D20","E20",
    "B21","C21","D21","E21"
)

foreach ($addr in $emptyCells) {
    $battingExtra.Range($addr).ClearContents()
}

# ---------------------------------------------------------------------------
# Part 2: add the "ODI Bowling Extra" sheet at the end of the workbook
# ---------------------------------------------------------------------------
$bowlingExtra = $wb.Worksheets.Add()
$bowlingExtra.Name = "ODI Bowling Extra"

# Move it to be the last tab (right after "ODI Batting Extra")
$lastIndex = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($lastIndex)
$bowlingExtra.Move($null, $lastSheet)

# Copy the header styling (bold / bordered / centered) from the sibling sheet
$battingExtra.Range("A1:C1").Copy()
$bowlingExtra.Range("A1:C1").PasteSpecial(-4122)

$bowlingExtra.Cells.Item(1,1).Value = "MATCH_CODE"
$bowlingExtra.Cells.Item(1,2).Value = "MAIDEN_OVERS"
$bowlingExtra.Cells.Item(1,3).Value = "PERCENT_WICKETS_OF_ALL"

# MATCH_CODE, MAIDEN_OVERS, PERCENT_WICKETS_OF_ALL -- "" means leave blank
$rows = @(
    @("3892","0",""),
    @("3898","0","20.00%"),
    @("3900","",""),
    @("3905","",""),
    @("3909","0","10.00%"),
    @("3939","0",""),
    @("3943","0",""),
    @("3944","",""),
    @("4378","0",""),
    @("4379","0",""),
    @("4385","",""),
    @("4387","0","10.00%"),
    @("4394","0",""),
    @("4414","",""),
    @("4417","0","10.00%"),
    @("4449","0","10.00%"),
    @("4450","",""),
    @("4486","",""),
    @("4520","0","10.00%"),
    @("4522","","")
)

$r = 2
foreach ($row in $rows) {
    for ($c = 1; $c -le 3; $c++) {
        $val = $row[$c - 1]
        if ($val -ne "") {
            $cell = $bowlingExtra.Cells.Item($r, $c)
            $cell.NumberFormat = "@"
            $cell.Value = $val
        }
    }
    $r = $r + 1
}
